$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the bottom rule under the year header row (B3:J3) - the header box
# keeps only its top edge now.
$ws.Range("B3:J3").Borders.Item(9).LineStyle = -4142

# Add the new 2023 column (K) following the existing year columns (B..J = 2014..2022)
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 603
$ws.Range("K5").Value = 254
$ws.Range("K6").Value = 349

# Match formatting of the preceding column (J) for the new column
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122) # xlPasteFormats

# K is now the rightmost column of the table, so it gets a right border to
# close off the table outline.
$ws.Range("K3:K6").Borders.Item(10).LineStyle = 1
$ws.Range("K3:K6").Borders.Item(10).Color = 0

# Extend the custom column width (8.7109375) that previously covered B:J to
# also cover the new columns through O, matching the widened "left table".
$ws.Range("B1:O1").ColumnWidth = 8.7109375
